$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").Value = '0.0₃0972'
$ws.Range("E35").Value = '  +14.21%  '
$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").Value = '67.52'
$ws.Range("E36").Value = '  -3.08%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.143'
$ws.Range("E42").Value = '  -2.46%  '
$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D48").Value = '2.56'
$ws.Range("E48").Value = '  -8.79%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '3.33'
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("D2").Value = '72.006.26'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '3.891.32'
$ws.Range("E3").Value = '  -2.08%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '592.58'
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("D6").Value = '167.99'
$ws.Range("E6").Value = '  +10.61%  '
$ws.Range("D7").Value = '0.673'
$ws.Range("E7").Value = '  -0.76%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").Value = '0.767'
$ws.Range("E9").Value = '  +2.46%  '
$ws.Range("D10").Value = '0.183'
$ws.Range("E10").Value = '  +9.02%  '
$ws.Range("D11").Value = '54.53'
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("D12").Value = '0.0000324'
$ws.Range("E12").Value = '  +2.22%  '
$ws.Range("D13").Value = '11.32'
$ws.Range("D14").Value = '4.524.21'
$ws.Range("E14").Value = '  -2.29%  '
$ws.Range("D15").Value = '3.914.34'
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").Value = '21.00'
$ws.Range("E16").Value = '  +2.65%  '
$ws.Range("D17").Value = '13.94'
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("D18").Value = '1.21'
$ws.Range("E18").Value = '  -4.00%  '
$ws.Range("D19").Value = '71.974.61'
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("D21").Value = '440.64'
$ws.Range("E21").Value = '  +1.64%  '
$ws.Range("D22").Value = '4.71'
$ws.Range("E22").Value = '  -0.48%  '
$ws.Range("D23").Value = '94.08'
$ws.Range("E23").Value = '  -1.88%  '
$ws.Range("D24").Value = '3.28'
$ws.Range("E24").Value = '  -5.97%  '
$ws.Range("D25").Value = '13.88'
$ws.Range("E25").Value = '  -2.44%  '
$ws.Range("D26").Value = '4.17'
$ws.Range("E26").Value = '  -4.95%  '
$ws.Range("D27").Value = '11.06'
$ws.Range("E27").Value = '  -4.38%  '
$ws.Range("D28").Value = '5.93'
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("D29").Value = '10.16'
$ws.Range("E29").Value = '  -4.37%  '
$ws.Range("D30").Value = '35.18'
$ws.Range("E30").Value = '  -2.90%  '
$ws.Range("D31").Value = '7.77'
$ws.Range("E31").Value = '  -2.70%  '
$ws.Range("D32").Value = '50.10'
$ws.Range("E32").Value = '  -0.86%  '
$ws.Range("D33").Value = '13.60'
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("D34").Value = '0.126'
$ws.Range("E34").Value = '  -4.59%  '
$ws.Range("D37").Value = '624.05'
$ws.Range("E37").Value = '  -8.28%  '
$ws.Range("D38").Value = '0.420'
$ws.Range("E38").Value = '  -4.57%  '
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("D40").Value = '3.31'
$ws.Range("E40").Value = '  +0.40%  '
$ws.Range("D43").Value = '3.20'
$ws.Range("E43").Value = '  +42.90%  '
$ws.Range("D44").Value = '0.0470'
$ws.Range("E44").Value = '  -3.40%  '
$ws.Range("D45").Value = '10.28'
$ws.Range("E45").Value = '  -7.14%  '
$ws.Range("D46").Value = '0.145'
$ws.Range("E46").Value = '  -2.83%  '
$ws.Range("D47").Value = '2.84'
$ws.Range("E47").Value = '  -15.74%  '
$ws.Range("D50").Value = '2.814.05'
$ws.Range("E50").Value = '  +1.32%  '
$ws.Range("D51").Value = '0.000273'
$ws.Range("E51").Value = '  +2.49%  '
